# #96 another one test for FUNCEXEC
# Adds a new test row (27) to Sheet1 exercising FUNCEXEC with a range argument
# (K27:L27), mirroring the pattern used by the existing FUNCEXEC test rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New inputs for the FUNCEXEC range-argument test.
$ws.Range("K27").Value = 3
$ws.Range("L27").Value = 7

# Expected-value column for the new test row.
$ws.Range("B27").Value = 10

# The actual FUNCEXEC formula under test.
$ws.Range("A27").Formula = '=FUNCEXEC("DEF_1", K27:L27)'

# Move the active selection below the newly added row, matching the
# author's final cursor position after inserting the row.
$ws.Range("B28").Select()
